$d = $word.ActiveDocument

# --- 1. Merge the three runs in the "Desc: This is ONLY for making Users..."
#        paragraph into a single run by replacing the text that spans the
#        run boundaries with itself (same formatting on both sides causes
#        the engine to coalesce adjacent runs).
$d.Content.Find.Execute("Users,(NOT HERE). The creators", $true, $false, $false, $false, $false, $true, 1, $false, "Users,(NOT HERE). The creators", 2) | Out-Null

# --- 2. Merge the two runs in the "Collecting: ..." paragraph into one.
$d.Content.Find.Execute("what they want to add to LearnR 5. First Name", $true, $false, $false, $false, $false, $true, 1, $false, "what they want to add to LearnR 5. First Name", 2) | Out-Null

# --- 3. Split the "Generate: ... Organization Goals(string) ..." run into
#        three runs by inserting "array of " before "string)". The engine
#        coalesces a freshly inserted run into its identically-formatted
#        neighbours, so we temporarily give the inserted text a different
#        size, then restore the matching size afterwards - by then the
#        run boundaries have already been committed as separate runs.
$rng = $d.Content
$rng.Find.Execute("string) 4. List of all Users under that organization", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertStart = $rng.Start
$insertPoint = $d.Range($insertStart, $insertStart)
$insertPoint.InsertBefore("array of ")
$insertedRange = $d.Range($insertStart, $insertStart + 9)
$insertedRange.Font.Size = 28

$fixupRange = $d.Content
$fixupRange.Find.Execute("array of ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fixupRange.Font.Size = 12

# --- 4. styles.xml Normal style pPr: add <w:suppressAutoHyphens w:val="true"/>
#        (Hyphenation = $false means "suppress automatic hyphenation").
$d.Styles("Normal").ParagraphFormat.Hyphenation = $false
